# Add new feature "proveedor facturas y productos": two new sheets (Hoja2, Hoja3)
# inserted into the workbook, plus updated selections on the pre-existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert "Hoja2" right after "MODULOS" (i.e. before "Hoja1").
# ---------------------------------------------------------------------------
$modulos = $wb.Worksheets.Item("MODULOS")
$hoja2 = $wb.Worksheets.Add($null, $modulos)
$hoja2.Name = "Hoja2"

$hoja2.Range("A1").Value = "1.- Definir campos proveedores"
$hoja2.Range("A2").Value = "2.- Los productos pueden estar definidos a 1 o mas proveedores?"
$hoja2.Range("A3").Value = "3.- Campos de los productos validar"
$hoja2.Range("B2").Value = "Todo"

$hoja2.Columns.Item(1).ColumnWidth = 58.5

$hoja2.Range("D6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert "Hoja3" right after "BASE DE DATOS" (i.e. at the very end).
# ---------------------------------------------------------------------------
$baseDeDatos = $wb.Worksheets.Item("BASE DE DATOS")
$hoja3 = $wb.Worksheets.Add($null, $baseDeDatos)
$hoja3.Name = "Hoja3"

$hoja3.Range("A1").Value = "nombre"
$hoja3.Range("A3").Value = "plazo"
$hoja3.Range("A4").Value = "cuota"
$hoja3.Range("B1").Value = "Video Vigilancia 24/7"
$hoja3.Range("B2").Value = "2 camaras"
$hoja3.Range("B3").Value = "18 meses"
$hoja3.Range("A5").Value = "activo"
$hoja3.Range("B5").Value = "si"
$hoja3.Range("A2").Value = "cantidad"
$hoja3.Range("B4").Value = 100

$hoja3.Range("B1:B5").HorizontalAlignment = -4108
$hoja3.Range("B1:B5").VerticalAlignment = -4108

$hoja3.Columns.Item(1).ColumnWidth = 36.5
$hoja3.Columns.Item(2).ColumnWidth = 33.5

$hoja3.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update the selection on the two pre-existing sheets.
# ---------------------------------------------------------------------------
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A8").Select() | Out-Null

$baseDeDatos.Range("C23").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore MODULOS as the active sheet with its new selection.
# ---------------------------------------------------------------------------
$modulos.Activate()
$modulos.Range("C24").Select() | Out-Null
